# Add a "break_trial" column (Z) to the trial list sheet.
# The header "break_trial" already exists in Z1 (shared string), so this
# script only needs to populate Z2:Z115 with the per-trial flag (0/1) and
# update the sheet's view/selection state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Per-row break_trial values for rows 2..115 (1-indexed list, index 0 -> row 2).
$breakTrialValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

# A handful of rows (51, 54, 94, 95) already have an empty, pre-styled Z
# cell (bold row highlighting carried over from neighbouring columns).
# Writing a value into them should drop that stray formatting, same as
# every other (previously blank / nonexistent) cell in the column.
$preStyledRows = @(51, 54, 94, 95)

for ($i = 0; $i -lt $breakTrialValues.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 26)   # column 26 = Z
    $cell.Value = $breakTrialValues[$i]
    if ($preStyledRows -contains $row) {
        $cell.ClearFormats()
    }
}

# Update the window/selection to mirror the edited author's view: scrolled
# over to column L, with the whole new Z column (Z1:Z115) selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$ws.Range("Z1:Z115").Select() | Out-Null

Write-Host "break_trial column populated"
